# PCC.ODH 1.0.0 for TI
# Applies the diff to the Metadata sheet: inserts an "Identifier" row,
# updates Version/Publisher text, expands the single "Contact" row into
# three Contact rows, and shifts everything else down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$xlPasteFormats = -4122

# --- Insert a new row for "Identifier" right after the URL row (old row 3) ---
$ws.Range("A3:B3").EntireRow.Insert()
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B3").PasteSpecial($xlPasteFormats)

$ws.Range("A3").Value = "Identifier"
$ws.Range("B3").Value = "OID:1.3.6.1.4.1.19376.1.5.3.1.3.43.42.8"

# Version text changes (now on row 4 after the insert above)
$ws.Range("B4").Value = "1.0.0"

# Publisher text change (now on row 10)
$ws.Range("B10").Value = "IHE Patient Care Coordination Committee"

# --- Expand the single "Contact" row (now row 11) into three Contact rows ---
$ws.Range("A12:B13").EntireRow.Insert()
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B13").PasteSpecial($xlPasteFormats)

$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "null (https://www.ihe.net/ihe_domains/patient_care_coordination/)"
$ws.Range("A12").Value = "Contact"
$ws.Range("B12").Value = "null (pcc@ihe.net)"
$ws.Range("A13").Value = "Contact"
$ws.Range("B13").Value = "IHE Patient Care Coordination Committee (pcc@ihe.net)"
